# edit.ps1 - apply the "remove files and commit" change to GitNote.docx
#
# Strategy: for every paragraph we need to change, locate it with
# Find.Execute (text search), grab its Paragraph object, then replace the
# *whole paragraph's* Range via Range.InsertXML with a fully-formed
# <w:p>...</w:p> fragment (or several, for multi-paragraph inserts). This
# reliably preserves/sets pPr (style, numPr, indent) and lets us create
# separate <w:r> runs and <w:proofErr/> markers exactly as authored by Word.

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyInner) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyInner + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) "Ignored files: Git just cannot see them" -> append
#    " (different with untracked)" as three separate runs.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Ignored files: Git just cannot see them", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph 'Ignored files...' not found" }
$para = $rng.Paragraphs(1)
$target = $para.Range
$body = '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr></w:pPr>' +
        '<w:r><w:t>Ignored files: Git just cannot see them</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
        '<w:r><w:t>different with untracked</w:t></w:r>' +
        '<w:r><w:t>)</w:t></w:r>' +
        '</w:p></w:body>'
$target.InsertXML((New-PkgXml $body))

# ---------------------------------------------------------------------
# 2) "Remove files" -> append ": git rm" as a second run.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Remove files", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph 'Remove files' not found" }
$para = $rng.Paragraphs(1)
$target = $para.Range
$body = '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="22"/></w:numPr></w:pPr>' +
        '<w:r><w:t>Remove files</w:t></w:r>' +
        '<w:r><w:t>: git rm</w:t></w:r>' +
        '</w:p></w:body>'
$target.InsertXML((New-PkgXml $body))

# Remember the paragraph right after "Remove files" -- in the original
# document it is a single empty "ListParagraph" paragraph; the diff
# replaces it with four new paragraphs.
$afterRemoveFiles = $para.Next()

$target = $afterRemoveFiles.Range
$body = '<w:body>' +
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr></w:pPr>' +
            '<w:r><w:t xml:space="preserve">Try with </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
            '<w:r><w:t>df.RData</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
            '<w:r><w:t xml:space="preserve">: not </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:t>gonna</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:t xml:space="preserve"> work, because it is ignored by Git</w:t></w:r>' +
        '</w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr></w:pPr>' +
            '<w:r><w:t>Git rm: remove from your device. Can' + [char]0x2019 + 't get it back. Use with caution!</w:t></w:r>' +
        '</w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr></w:pPr>' +
            '<w:r><w:t xml:space="preserve">Git rm: cached: move it back to ' + [char]0x201C + 'untracked' + [char]0x201D + '. Git still sees it. </w:t></w:r>' +
        '</w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1080"/></w:pPr>' +
            '<w:r><w:t>If use ' + [char]0x201C + '</w:t></w:r>' +
            '<w:proofErr w:type="gramStart"/>' +
            '<w:r><w:t>add .</w:t></w:r>' +
            '<w:proofErr w:type="gramEnd"/>' +
            '<w:r><w:t>' + [char]0x201D + ' it will be add back to git</w:t></w:r>' +
        '</w:p>' +
        '</w:body>'
$target.InsertXML((New-PkgXml $body))

# ---------------------------------------------------------------------
# 3) Move <w:lastRenderedPageBreak/> from the "Undo changes:" run to the
#    "But mostly can be done on GitHub, so omitted" run.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("But mostly can be done on GitHub, so omitted", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph 'But mostly...' not found" }
$para = $rng.Paragraphs(1)
$target = $para.Range
$body = '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
        '<w:r><w:lastRenderedPageBreak/><w:t>But mostly can be done on GitHub, so omitted</w:t></w:r>' +
        '</w:p></w:body>'
$target.InsertXML((New-PkgXml $body))

$rng = $d.Content
$found = $rng.Find.Execute("Undo changes:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph 'Undo changes:' not found" }
$para = $rng.Paragraphs(1)
$target = $para.Range
$body = '<w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr>' +
        '<w:r><w:t>Undo changes:</w:t></w:r>' +
        '</w:p></w:body>'
$target.InsertXML((New-PkgXml $body))

Write-Host "edits applied"
